$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.315.43"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "1.585.28"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'209.49"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "'0.504"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("D7").Value = "'1.01"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").Value = "'19.48"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").Value = "'0.0842"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "1.808.23"
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").Value = "1.584.21"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "'0.516"
$ws.Range("E15").Value = "  -1.48%  "
$ws.Range("D16").Value = "'64.25"
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").Value = "26.312.02"
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").Value = "'7.22"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").Value = "'1.01"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "'206.89"
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("D22").Value = "'4.27"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").Value = "'2.21"
$ws.Range("E23").Value = "  -3.77%  "
$ws.Range("D24").Value = "'8.81"
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("D25").Value = "'144.44"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("D27").Value = "'7.02"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").Value = "'0.113"
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("D29").Value = "'15.28"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").Value = "'0.0504"
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("D32").Value = "'3.23"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("D33").Value = "'2.94"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("D34").Value = "'1.30"
$ws.Range("E34").Value = "  +13.75%  "
$ws.Range("D35").Value = "1.284.96"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D37").Value = "'0.605"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").Value = "'0.0167"
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("D40").Value = "'0.817"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("D41").Value = "'5.44"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("E42").Value = "  -1.71%  "
$ws.Range("E43").Value = "  -4.83%  "
$ws.Range("D44").Value = "'62.21"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("D45").Value = "1.720.21"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").Value = "'88.71"
$ws.Range("E46").Value = "  -2.84%  "
$ws.Range("D47").Value = "'1.56"
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "'7.44"
$ws.Range("E51").Value = "  +0.22%  "
